$wb = $excel.ActiveWorkbook

# ---------- Sheet: "Test Cases" ----------
# Add a new test case row (TC_002) right after the existing TC_001 row,
# re-using the same row formatting as the row above it.
$ws1 = $wb.Worksheets.Item("Test Cases")
$ws1.Range("A2:C2").Copy($ws1.Range("A3:C3"))
$ws1.Range("A3").Value = "TC_002_Validate_NearBy_Homes_page"
$ws1.Range("B3").Value = "Validate the near by homes page count"
$ws1.Range("C3").Value = "Y"
$ws1.Range("A3").Select()

# ---------- Sheet: "Test Steps" ----------
# Duplicate the first three steps of TC_001 (cookie accept, skip
# registration, verify search title) as the first three steps of the new
# TC_002 test case, then append one brand-new step specific to TC_002.
$ws2 = $wb.Worksheets.Item("Test Steps")

$ws2.Range("A2:F4").Copy($ws2.Range("A12:F14"))
$ws2.Range("A12").Value = "TC_002_Validate_NearBy_Homes_page"
$ws2.Range("A13").Value = "TC_002_Validate_NearBy_Homes_page"
$ws2.Range("A14").Value = "TC_002_Validate_NearBy_Homes_page"

$ws2.Range("A2:F2").Copy($ws2.Range("A15:F15"))
$ws2.Range("A15").Value = "TC_002_Validate_NearBy_Homes_page"
$ws2.Range("C15").Value = "scroll_Down"
$ws2.Range("D15").Value = "nearyoutext"
$ws2.Range("B15").Value = "NearBy homes location selector"

$ws2.Range("C21").Select()
